# "Generate Report for handoff"
# Updates the zh-cn and de-de handoff-status sheets: the handoff package is now
# ready, so status moves from "Handoff transform failed" / "Ignored" to
# "Ready for handoff" / "Include", and the newly produced .xlf handoff file
# plus its timestamp are recorded.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/8618ad2fdf283a6d8e0cd2a7216d539d46093725"

# The Overview sheet mirrors the per-locale status for each source file
# (the "Handoff transform failed" text is the very same shared string used
# on the zh-cn/de-de detail sheets), so it also needs to reflect the new
# "Ready for handoff" status for the a6d5f17a-...md row.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"

$langs = @(
    @{ Sheet = "zh-cn"; File = "a6d5f17a-6127-422a-a9c7-e2c1c8202ce2.6a84fe65a7cfe9ef7597db12c5983f89d006dece.zh-cn.xlf"; Stamp = "2016-01-20 08:11:20" },
    @{ Sheet = "de-de"; File = "a6d5f17a-6127-422a-a9c7-e2c1c8202ce2.6a84fe65a7cfe9ef7597db12c5983f89d006dece.de-de.xlf"; Stamp = "2016-01-20 08:11:32" }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Row 2 ("a6d5f17a-...md" source file): handoff succeeded.
    $ws.Range("B2").Value = "Ready for handoff"
    $ws.Range("H2").Value = "Include"

    # Record the newly generated handoff (.xlf) file and when it was produced.
    $ws.Hyperlinks.Add($ws.Range("C2"), "$repoBase/$($lang.File)", "", "", $lang.File)
    $ws.Range("D2").Value = $lang.Stamp
}
